$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a serial date value that was incremented by
# one day (45177 -> 45178) for every data row (rows 2 through 140).
$ws.Range("C2:C140").Value = 45178
